$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay text (matching the source
# inline-string cells): force a Text number format before assigning, then drop
# back to the Normal style so the cell does not keep a stray style index.
$textForceCells = @('D5', 'D6', 'D10', 'D11', 'D16', 'D18', 'D20', 'D22', 'D24', 'D29', 'D31', 'D32', 'D39', 'D40', 'D44', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.881.48'

$ws.Range('D3').Value = '2.267.24'
$ws.Range('E3').Value = '  +2.22%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '301.84'
$ws.Range('E5').Value = '  +3.53%  '

$ws.Range('D6').Value = '92.15'
$ws.Range('E6').Value = '  +6.24%  '

$ws.Range('E7').Value = '  +3.43%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  +4.30%  '

$ws.Range('D10').Value = '54.34'
$ws.Range('E10').Value = '  +7.85%  '

$ws.Range('D11').Value = '32.27'
$ws.Range('E11').Value = '  +6.24%  '

$ws.Range('E12').Value = '  +2.56%  '

$ws.Range('E13').Value = '  +2.11%  '

$ws.Range('E14').Value = '  +3.80%  '

$ws.Range('D15').Value = '2.618.04'
$ws.Range('E15').Value = '  +2.18%  '

$ws.Range('D16').Value = '14.16'
$ws.Range('E16').Value = '  +2.86%  '

$ws.Range('D17').Value = '2.269.77'
$ws.Range('E17').Value = '  +0.53%  '

$ws.Range('D18').Value = '0.758'
$ws.Range('E18').Value = '  +3.66%  '

$ws.Range('D19').Value = '41.780.52'
$ws.Range('E19').Value = '  +4.78%  '

$ws.Range('D20').Value = '12.13'
$ws.Range('E20').Value = '  +9.27%  '

$ws.Range('E21').Value = '  +2.18%  '

$ws.Range('D22').Value = '5.95'
$ws.Range('E22').Value = '  +3.65%  '

$ws.Range('E23').Value = '  +2.29%  '

$ws.Range('D24').Value = '241.67'
$ws.Range('E24').Value = '  +2.01%  '

$ws.Range('E25').Value = '  +4.40%  '

$ws.Range('E27').Value = '  +3.91%  '

$ws.Range('E28').Value = '  +2.71%  '

$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  +4.51%  '

$ws.Range('E30').Value = '  -11.96%  '

$ws.Range('D31').Value = '159.78'
$ws.Range('E31').Value = '  +1.39%  '

$ws.Range('D32').Value = '33.81'
$ws.Range('E32').Value = '  +6.52%  '

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('E34').Value = '  +4.18%  '

$ws.Range('E35').Value = '  +4.58%  '

$ws.Range('E36').Value = '  +2.98%  '

$ws.Range('E37').Value = '  +2.04%  '

$ws.Range('E38').Value = '  +5.61%  '

$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '0.116'
$ws.Range('E39').Value = '  +3.55%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '16.53'
$ws.Range('E40').Value = '  +9.03%  '

$ws.Range('E41').Value = '  +5.10%  '

$ws.Range('E42').Value = '  +5.96%  '

$ws.Range('D43').Value = '2.070.95'
$ws.Range('E43').Value = '  -0.62%  '

$ws.Range('D44').Value = '19.81'
$ws.Range('E44').Value = '  +10.56%  '

$ws.Range('E45').Value = '  +3.31%  '

$ws.Range('D46').Value = '10.19'
$ws.Range('E46').Value = '  +4.21%  '

$ws.Range('D47').Value = '2.92'
$ws.Range('E47').Value = '  +8.52%  '

$ws.Range('E48').Value = '  +2.39%  '

$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').Value = '1.15'
$ws.Range('E49').Value = '  +3.79%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.51'
$ws.Range('E50').Value = '  +3.74%  '

$ws.Range('D51').Value = '51.71'
$ws.Range('E51').Value = '  +5.83%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
